# Regenerate S-val data (filter save games) for rows 2-48, columns B:G
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$n_cols = 6
# New values per row: TB, d2S, K, IP, Win, sum
$rows = @(
    @(1.445647641019636, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 1, 4.327115817150455),
    @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 1, 6.15379541431027),
    @(3.272327238179451, 1.626987699542094, 18.71679738969934, 0.5333859586016987, 1, 24.14949828602258),
    @(0.1169995834814548, 0.3048912486333797, 0.1496068669990043, 0.5333859586016987, 1, 1.104883657715537),
    @(0.04172184405617529, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 1, 2.923190020186994),
    @(3.272327238179451, 1.626987699542094, 18.71679738969934, 0.5333859586016987, 1, 24.14949828602258),
    @(0.2881169905109251, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 2.598097515653722),
    @(0.2881169905109251, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 0, 5.671859677732939),
    @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 1, 8.656069925401464),
    @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 5.582307763322248),
    @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 3.755628166162433),
    @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 1, 6.15379541431027),
    @(1.445647641019636, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 0, 4.327115817150455),
    @(1.445647641019636, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 1, 4.327115817150455),
    @(0.1169995834814548, 0.3048912486333797, 3.223369029078222, 0.5333859586016987, 0, 4.178645819794754),
    @(1.445647641019636, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 1, 6.82939032824165),
    @(3.272327238179451, 1.626987699542094, 18.71679738969934, 0.5333859586016987, 1, 24.14949828602258),
    @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 5.582307763322248),
    @(0.6545652718822623, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 1, 3.536033448013082),
    @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 0, 3.755628166162433),
    @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 0, 6.15379541431027),
    @(0.6545652718822623, 0.3048912486333797, 18.71679738969934, 0.5333859586016987, 1, 20.20963986881668),
    @(0.1169995834814548, 0.3048912486333797, 18.71679738969934, 0.5333859586016987, 0, 19.67207418041587),
    @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 3.755628166162433),
    @(0.6545652718822623, 0.3048912486333797, 3.223369029078222, 0.5333859586016987, 0, 4.716211508195562),
    @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 1, 6.15379541431027),
    @(3.272327238179451, 1.626987699542094, 18.71679738969934, 0.5333859586016987, 1, 24.14949828602258),
    @(0.6545652718822623, 0.3048912486333797, 0.7210945179870265, 0.5333859586016987, 1, 2.213936997104367),
    @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 5.582307763322248),
    @(1.445647641019636, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 1, 6.82939032824165),
    @(3.272327238179451, 1.626987699542094, 18.71679738969934, 0.5333859586016987, 1, 24.14949828602258),
    @(0.04172184405617529, 0.04103571897497393, 3.223369029078222, 13.86384647080068, 0, 17.16997306291006),
    @(3.272327238179451, 1.626987699542094, 18.71679738969934, 0.5333859586016987, 1, 24.14949828602258),
    @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 1, 6.15379541431027),
    @(3.272327238179451, 1.626987699542094, 18.71679738969934, 0.5333859586016987, 0, 24.14949828602258),
    @(3.272327238179451, 9.983522426115931, 0.7210945179870265, 13.86384647080068, 1, 27.84079065308309),
    @(3.272327238179451, 1.626987699542094, 18.71679738969934, 0.5333859586016987, 1, 24.14949828602258),
    @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 1, 8.656069925401464),
    @(3.272327238179451, 1.626987699542094, 189.6080260415259, 0.5333859586016987, 1, 195.0407269378492),
    @(0.04172184405617529, 0.002658071450198252, 0.7210945179870265, 0.5333859586016987, 1, 1.298860392095099),
    @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 0, 6.15379541431027),
    @(1.445647641019636, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 1, 4.327115817150455),
    @(0.2881169905109251, 0.3048912486333797, 0.7210945179870265, 0.5333859586016987, 1, 1.84748871573303),
    @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 1, 6.15379541431027),
    @(3.272327238179451, 1.626987699542094, 18.71679738969934, 0.5333859586016987, 0, 24.14949828602258),
    @(3.272327238179451, 1.626987699542094, 18.71679738969934, 0.5333859586016987, 1, 24.14949828602258),
    @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 0, 8.656069925401464)
)

$data = New-Object "object[,]" 47,6
for ($i = 0; $i -lt $rows.Count; $i++) {
    for ($j = 0; $j -lt $n_cols; $j++) {
        $data[$i,$j] = $rows[$i][$j]
    }
}

$ws.Range("B2:G48").Value = $data
